# Applies the "Genes in the Genomic Era" -> "The Vital Role of Chemistry in
# Everyday Life" rewrite described by the target diff.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $range = $d.Content
    $ok = $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Output ("MISSING: " + $old)
    }
}

# --- Title ---------------------------------------------------------------
Replace-Text "Genes in the Genomic Era" "The Vital Role of Chemistry in Everyday Life: An Exploration for High School Students"

# --- Author name -----------------------------------------------------------
Replace-Text " Ashley Isabella" " Jeremy Peterson"

# --- Author email (split across several runs) ------------------------------
Replace-Text "ashleyisabella@biomed" "jeremy"
Replace-Text "edu" "peterson@successacademy"

# Append the ".org" part as new text after the email we just built.
$emailPara = $d.Paragraphs.Item(3)
$endPos = $emailPara.Range.End
$insertPoint = $d.Range($endPos - 1, $endPos - 1)
$insertPoint.InsertAfter(".")
$endPos = $emailPara.Range.End
$insertPoint = $d.Range($endPos - 1, $endPos - 1)
$insertPoint.InsertAfter("org")

# --- Body paragraph: intro sentences ---------------------------------------
Replace-Text "Genomics, the study of genomes, has revolutionized our understanding of life" "Chemistry, the science that focuses on the composition, structure, properties, and changes of matter, plays a pivotal role in our daily lives, shaping everything from the food we eat to the medicines we take"
Replace-Text " The sequencing of the human genome in 2003 marked a watershed moment, providing an unprecedented roadmap of human DNA" " This interdisciplinary subject unlocks the mysteries of the chemical reactions that govern the world around us"
Replace-Text " Through this lens of investigation, genomic era introduced advances in comprehending heritable traits, the intricacies of evolution, and applications in medical diagnostics and treatment" " To truly understand how chemistry impacts our lives, it is essential to delve into its diverse applications and the fundamental principles that underpin them"
Replace-Text " Genomics, in its entirety, is radically transforming biology and medicine, enabling us to peer into the molecular underpinnings of life in ways previously unfathomable" " In this essay, we will embark on an enlightening journey to explore the vital role of chemistry in myriad aspects of everyday life, encompassing food and agriculture, medicine, materials science, and energy sources"

# --- Body paragraph: food & agriculture section -----------------------------
Replace-Text "In this genomics era, scientists have been able to identify genetic variations associated with a range of diseases, leading to the development of personalized medicine and targeted therapies" "From the nutritional components of the food we consume to the intricate processes that sustain agricultural ecosystems, chemistry lies at the heart of our sustenance"
Replace-Text " Genetic testing has become increasingly accessible, allowing individuals to understand their genetic predispositions and make informed decisions about their health" " Chemical reactions and interactions orchestrate the growth of crops and livestock, ensuring food production can meet the growing demands of a burgeoning population"
Replace-Text " The ability to sequence and analyze genomic data has also fueled advancements in fields such as tracing evolutionary history, underpinnings of human behavior" " The study of chemistry empowers us to comprehend these processes, develop sustainable agricultural practices, and address global food security challenges"

# --- Body paragraph: medicine section ---------------------------------------
Replace-Text "The profound implications of genomics extend beyond health and scientific research" "The realm of medicine is profoundly intertwined with chemistry"
Replace-Text " The insights gleaned from studying genomes have significant societal, ethical, and legal dimensions" " Pharmaceutical research harnesses chemistry to design and synthesize groundbreaking drugs, alleviating diseases and mitigating debilitating conditions"
Replace-Text " The accessibility and privacy of genetic information, the potential for genetic discrimination, and the ethical considerations surrounding germline editing all demand careful consideration" " The marvels of vaccines and antibiotics, along with cutting-edge cancer treatments, underscore the transformative impact of chemistry on human health and well-being"
Replace-Text " Genomics has sparked discussions on the boundaries of human enhancement, the nature of identity, and the very meaning of life in this era of unprecedented genetic knowledge" " Understanding medicinal chemistry enables us to make informed decisions about our health, appreciate the intricacies of drug development, and marvel at the life-saving discoveries that stem from scientific advancements"

# --- Body paragraph: new Materials science & Energy sections ---------------
# Inserted right after "...scientific advancements" and before the sentence's
# closing period (which is an untouched run carried over from the original).
$vt = [string][char]11
$materialsAndEnergy = (
    $vt + $vt +
    "Chemistry drives the development of innovative materials that shape our modern world." +
    " From the polymers in synthetic fabrics to the alloys in infrastructure, chemistry empowers us to tailor materials with specific properties and applications." +
    " Advances in materials chemistry have led to groundbreaking technologies, including lightweight composites, energy-efficient electronics, and sustainable construction materials." +
    " Comprehending the chemical principles behind materials science equips us to appreciate the engineering feats that underpin modern society and envision the possibilities of future material innovations." +
    $vt + $vt +
    "Lastly, our quest for sustainable energy sources hinges on the study of chemistry." +
    " The transition to cleaner and renewable energy alternatives, such as solar cells and fuel cells, relies heavily on chemical processes and reactions." +
    " Understanding electrochemistry and energy storage technologies gives us the tools to address pressing environmental issues, mitigate climate change, and secure a brighter energy future for generations to come"
)

$anchor = $d.Content
$found = $anchor.Find.Execute("scientific advancements")
if (-not $found) {
    Write-Output "MISSING: scientific advancements anchor"
}
$anchor.Collapse(0)
$anchor.InsertAfter($materialsAndEnergy)

# --- Summary paragraph: final sentence (merging two runs into one) ---------
$oldSummaryTail = " As we navigate this era of unprecedented genetic knowledge, we are poised to unlock further breakthroughs in medicine, while thoughtfully addressing the complex implications and uncertainties that accompany this newfound understanding of life's molecular blueprint"
$newSummaryTail = " Chemistry serves as a cornerstone of scientific discovery and technological progress, inspiring us to continuously explore the wonders of the natural world and harness its potential for the betterment of society"
Replace-Text $oldSummaryTail $newSummaryTail

# --- New trailing empty paragraph at the very end of the document ----------
$lastPara = $d.Paragraphs.Last
$endOfDoc = $lastPara.Range.End
$endRange = $d.Range($endOfDoc, $endOfDoc)
$endRange.InsertParagraphAfter()

Write-Output "done"
